# Updated loading_percent values for case with 380 kV
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$columns = @('B','C','D','F','G','H','I','J','K')

$values = @(
    @(10.51428138743031, 6.295178173772804, 9.749605101042016, 37.22636074956853, 40.75682497576504, 17.24912843580258, 24.98568356910541, 11.42912107208592, 10.99743220493115),
    @(10.25162481761651, 6.09691674957317, 9.666224074990254, 37.22862395910788, 40.76062850409306, 17.29467743643058, 25.06508622881022, 11.41286008246627, 10.8230573012549),
    @(10.0891008916001, 5.973325603963222, 9.616668020548589, 37.23940729583802, 40.77516573233449, 17.3256819721062, 25.11882061272538, 11.405194638054, 10.71683852968686),
    @(10.02266018200231, 5.922578241031853, 9.596905256263529, 37.24616121230022, 40.78415073341661, 17.33907906616242, 25.1419671575487, 11.40265605814162, 10.67382258841155),
    @(10.01161813682646, 5.914131125360838, 9.593650287532421, 37.2474251085367, 40.78582730358489, 17.34134965867229, 25.14588598860177, 11.40226992203364, 10.66669770624542),
    @(10.08820556236248, 5.972642642319294, 9.616399720055814, 37.23948883197257, 40.77527452487947, 17.32585956441487, 25.11912772078896, 11.40515803023505, 10.71625723847362),
    @(10.4240380967561, 6.227251895521203, 9.720526694396737, 37.22519035266829, 40.75560022147625, 17.2642026087277, 25.01202597721648, 11.42303435887835, 10.93716514991149),
    @(11.06835404434612, 6.708408211974326, 9.936852019413783, 37.27173820561299, 40.81406144639858, 17.16744746967084, 24.84165788856678, 11.47638114389322, 11.37465719276499),
    @(11.52764846785426, 7.046764596021878, 10.10196054367815, 37.35139091302796, 40.91633543506256, 17.11115910623652, 24.74084983425159, 11.52655583040691, 11.69548772823308),
    @(11.73251311434905, 7.196663000610084, 10.17814468920389, 37.39746142060319, 40.97572798611214, 17.08877979185606, 24.7003176226241, 11.55172077894041, 11.84065684694228),
    @(11.80942803450829, 7.252793996891175, 10.20712539263505, 37.41631576689966, 41.00006313654852, 17.08077040390437, 24.68573800918148, 11.56158206366574, 11.89546433138668),
    @(11.79289367852285, 7.240734101087702, 10.20087840486929, 37.41219261736062, 40.99474021309993, 17.08247466184769, 24.68884373593784, 11.5594435783795, 11.88366871430928),
    @(11.73885465754016, 7.201293891770284, 10.18052644843439, 37.39898439163636, 40.97769313007247, 17.08811152531003, 24.6991027180584, 11.55252544986508, 11.84516949263633),
    @(11.70566567888236, 7.177051780988839, 10.16807671220251, 37.39107718828425, 40.96749128786722, 17.09162488058709, 24.70548688784007, 11.5483309696277, 11.82156459522275),
    @(11.51417139720195, 7.036882710326323, 10.09700138732547, 37.34857755060846, 40.91271247786893, 17.11268667747736, 24.74360615572948, 11.52495794473742, 11.68598034566746),
    @(11.39559608253768, 6.949823225164928, 10.05365792500772, 37.325020921925, 40.88239978790974, 17.12643467813439, 24.76835748536478, 11.51121556061711, 11.60256609632919),
    @(11.32701577417534, 6.899373222907999, 10.0288305366145, 37.3123980517141, 40.86617640238023, 17.13464579449248, 24.78309490682401, 11.50353188430613, 11.5545182134445),
    @(11.30373304586905, 6.882228934711127, 10.02044274193981, 37.30828340976517, 40.86089168579334, 17.13747804852971, 24.78817071178617, 11.50096834329983, 11.5382396156127),
    @(11.40825842605887, 6.959130147867876, 10.05826145735774, 37.32743273687871, 40.88550125633142, 17.12493974649109, 24.7656707773876, 11.51265566411831, 11.61145334289158),
    @(11.75474580442449, 7.212895998498981, 10.1865009299655, 37.40282579850825, 40.98265027223903, 17.0864432079962, 24.6960685093674, 11.55454850437914, 11.85648253458916),
    @(11.97730009746611, 7.375040244519399, 10.27106808336367, 37.46030511070523, 41.05688930433591, 17.06399538470462, 24.65506346073354, 11.58385980760603, 12.01563902590521),
    @(11.85889906886309, 7.288856294430975, 10.22587157538776, 37.4288788977813, 41.01628581989304, 17.07572769395496, 24.67653730400573, 11.56804068937028, 11.93080121415042),
    @(11.4025350509581, 6.954923724710363, 10.05617991409127, 37.32633948860445, 40.88409533221572, 17.12561464819142, 24.76688385723658, 11.51200391703828, 11.60743570095585),
    @(10.89615523524031, 6.580631391313382, 9.877155674622248, 37.25115621923936, 40.78783061161331, 17.19102871661758, 24.88348191551614, 11.4600070203651, 11.25618042752149)
)

$startRow = 2
for ($i = 0; $i -lt $values.Count; $i++) {
    $row = $startRow + $i
    $rowValues = $values[$i]
    for ($j = 0; $j -lt $columns.Count; $j++) {
        $col = $columns[$j]
        $ws.Range("$col$row").Value = $rowValues[$j]
    }
}

$wb.Save()